# Rewrite the small roster table into a Role/Hours table.
#
# Strategy: first overwrite the four existing name cells (A2:A5) in place so
# the workbook's shared-string table gets the new text values while
# preserving each cell's original shared-string slot/order. Then insert a
# new header row, fill in the final header + data values, and drop the
# now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: seed the new strings into the existing cells (keeps shared-string
# index assignment in the same order the final file expects: Accounting(0),
# Role(1), Hours(2), Developer(3)).
$ws.Range("A2").Value = "Accounting"
$ws.Range("A3").Value = "Role"
$ws.Range("A4").Value = "Hours"
$ws.Range("A5").Value = "Developer"

# Step 2: insert a new header row above the table.
$ws.Rows("1:1").Insert()

# Step 3: write the final header row and make it bold.
$ws.Range("A1").Value = "Role"
$ws.Range("B1").Value = "Hours"
$ws.Range("A1:B1").Font.Bold = $true

# Step 4: write the final two data rows.
$ws.Range("A2").Value = "Accounting"
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = "Developer"
$ws.Range("B3").Value = 9

# Step 5: remove the now-stale leftover rows (old Role/Hours/Developer rows
# that shifted down to 4:6 after the insert).
$ws.Rows("4:6").Delete()

# Match the final selection state.
$ws.Range("A4").Select()
